$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 296.66666
$ws.Range("I2").Value = 251.33333
$ws.Range("J2").Value = 432.66666
$ws.Range("K2").Value = 251.33333
$ws.Range("L2").Value = 432.66666
$ws.Range("M2").Value = -138.33333
$ws.Range("N2").Value = -658.66666
# Row 40
$ws.Range("H40").Value = 6718.75
$ws.Range("J40").Value = 7083.3335
$ws.Range("L40").Value = 7083.3335
$ws.Range("N40").Value = -7433.3335
# Row 51
$ws.Range("H51").Value = 7653.9414
$ws.Range("I51").Value = 7567.737
$ws.Range("J51").Value = 7763.1333
$ws.Range("K51").Value = 7567.737
$ws.Range("L51").Value = 7763.1333
$ws.Range("M51").Value = -7083.737
$ws.Range("N51").Value = -8731.133300000001
# Row 64
$ws.Range("H64").Value = 7771.3335
$ws.Range("I64").Value = 5540.143
$ws.Range("J64").Value = 9191.182000000001
$ws.Range("K64").Value = 5540.143
$ws.Range("L64").Value = 9191.182000000001
$ws.Range("M64").Value = -5292.143
$ws.Range("N64").Value = -9687.182000000001
# Row 67
$ws.Range("H67").Value = 7771.3335
$ws.Range("I67").Value = 5540.143
$ws.Range("J67").Value = 9191.182000000001
$ws.Range("K67").Value = 5540.143
$ws.Range("L67").Value = 9191.182000000001
$ws.Range("M67").Value = -4682.143
$ws.Range("N67").Value = -10907.182
# Row 98
$ws.Range("H98").Value = 254960.42
$ws.Range("I98").Value = 940.8946999999999
$ws.Range("K98").Value = 940.8946999999999
$ws.Range("M98").Value = 557.1053000000001
# Row 122
$ws.Range("H122").Value = 254960.42
$ws.Range("I122").Value = 940.8946999999999
$ws.Range("K122").Value = 2822.6841
$ws.Range("M122").Value = -372.6840999999999
# Row 135
$ws.Range("H135").Value = 1377.5625
$ws.Range("I135").Value = 932.2143
$ws.Range("K135").Value = 8389.9287
$ws.Range("M135").Value = -5854.9287

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3113.7402
$ws.Range("I32").Value = 2181.2837
$ws.Range("K32").Value = 2181.2837
$ws.Range("M32").Value = -1894.2837
# Row 74
$ws.Range("H74").Value = 10419674
$ws.Range("I74").Value = 12347767
$ws.Range("J74").Value = 7971
$ws.Range("K74").Value = 12347767
$ws.Range("L74").Value = 7971
$ws.Range("M74").Value = -12346893
$ws.Range("N74").Value = -9719
# Row 77
$ws.Range("H77").Value = 10419674
$ws.Range("I77").Value = 12347767
$ws.Range("J77").Value = 7971
$ws.Range("K77").Value = 61738835
$ws.Range("L77").Value = 39855
$ws.Range("M77").Value = -61734467
$ws.Range("N77").Value = -48591
# Row 97
$ws.Range("H97").Value = 740.86664
$ws.Range("I97").Value = 654.8182
$ws.Range("J97").Value = 977.5
$ws.Range("K97").Value = 654.8182
$ws.Range("L97").Value = 977.5
$ws.Range("M97").Value = -158.8182
$ws.Range("N97").Value = -1969.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# Row 80
$ws.Range("H80").Value = 532.2778
$ws.Range("I80").Value = 440.6
$ws.Range("K80").Value = 440.6
$ws.Range("M80").Value = 557.4
# Row 83
$ws.Range("H83").Value = 532.2778
$ws.Range("I83").Value = 440.6
$ws.Range("K83").Value = 2203
$ws.Range("M83").Value = 2789
# Row 86
$ws.Range("H86").Value = 5004.5
$ws.Range("I86").Value = 3877
$ws.Range("J86").Value = 7635.3335
$ws.Range("K86").Value = 3877
$ws.Range("L86").Value = 7635.3335
$ws.Range("M86").Value = -2754
$ws.Range("N86").Value = -9881.333500000001
# Row 89
$ws.Range("H89").Value = 5004.5
$ws.Range("I89").Value = 3877
$ws.Range("J89").Value = 7635.3335
$ws.Range("K89").Value = 19385
$ws.Range("L89").Value = 38176.6675
$ws.Range("M89").Value = -13769
$ws.Range("N89").Value = -49408.6675
# Row 94
$ws.Range("H94").Value = 1769.5385
$ws.Range("I94").Value = 1143.3684
$ws.Range("K94").Value = 1143.3684
$ws.Range("M94").Value = -692.3684000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2635.2307
$ws.Range("I16").Value = 1825.8
$ws.Range("J16").Value = 5333.3335
$ws.Range("K16").Value = 1825.8
$ws.Range("L16").Value = 5333.3335
$ws.Range("M16").Value = -1538.8
$ws.Range("N16").Value = -5907.3335
# Row 31
$ws.Range("H31").Value = 23662.434
$ws.Range("I31").Value = 2605.7297
$ws.Range("J31").Value = 72356.06
$ws.Range("K31").Value = 2605.7297
$ws.Range("L31").Value = 72356.06
$ws.Range("M31").Value = -2310.7297
$ws.Range("N31").Value = -72946.06
# Row 34
$ws.Range("H34").Value = 23662.434
$ws.Range("I34").Value = 2605.7297
$ws.Range("J34").Value = 72356.06
$ws.Range("K34").Value = 2605.7297
$ws.Range("L34").Value = 72356.06
$ws.Range("M34").Value = -2403.7297
$ws.Range("N34").Value = -72760.06
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
# Row 113
$ws.Range("H113").Value = 2635.2307
$ws.Range("I113").Value = 1825.8
$ws.Range("J113").Value = 5333.3335
$ws.Range("K113").Value = 1825.8
$ws.Range("L113").Value = 5333.3335
$ws.Range("M113").Value = 344.2
$ws.Range("N113").Value = -9673.333500000001
# Row 132
$ws.Range("H132").Value = 2636.7334
$ws.Range("I132").Value = 1984.2222
$ws.Range("K132").Value = 5952.6666
$ws.Range("M132").Value = -3422.6666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2281.7666
$ws.Range("I113").Value = 1590.4546
$ws.Range("J113").Value = 4182.875
$ws.Range("K113").Value = 1590.4546
$ws.Range("L113").Value = 4182.875
$ws.Range("M113").Value = 579.5454
$ws.Range("N113").Value = -8522.875
# Row 122
$ws.Range("H122").Value = 4135.75
$ws.Range("I122").Value = 3339.6875
$ws.Range("J122").Value = 5727.875
$ws.Range("K122").Value = 10019.0625
$ws.Range("L122").Value = 17183.625
$ws.Range("M122").Value = -7569.0625
$ws.Range("N122").Value = -22083.625
# Row 126
$ws.Range("H126").Value = 3614.8276
$ws.Range("I126").Value = 2488.8667
$ws.Range("K126").Value = 7466.6001
$ws.Range("M126").Value = -4996.6001
# Row 132
$ws.Range("H132").Value = 4625.846
$ws.Range("I132").Value = 4455.325
$ws.Range("J132").Value = 5194.25
$ws.Range("K132").Value = 13365.975
$ws.Range("L132").Value = 15582.75
$ws.Range("M132").Value = -10835.975
$ws.Range("N132").Value = -20642.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5979.364
$ws.Range("I7").Value = 4915.3955
$ws.Range("K7").Value = 4915.3955
$ws.Range("M7").Value = -4803.3955
# Row 22
$ws.Range("H22").Value = 4876.773
$ws.Range("I22").Value = 2532.4666
$ws.Range("J22").Value = 9900.286
$ws.Range("K22").Value = 2532.4666
$ws.Range("L22").Value = 9900.286
$ws.Range("M22").Value = -2237.4666
$ws.Range("N22").Value = -10490.286
# Row 27
$ws.Range("H27").Value = 4876.773
$ws.Range("I27").Value = 2532.4666
$ws.Range("J27").Value = 9900.286
$ws.Range("K27").Value = 2532.4666
$ws.Range("L27").Value = 9900.286
$ws.Range("M27").Value = -2425.4666
$ws.Range("N27").Value = -10114.286
# Row 126
$ws.Range("H126").Value = 5979.364
$ws.Range("I126").Value = 4915.3955
$ws.Range("K126").Value = 14746.1865
$ws.Range("M126").Value = -12276.1865
# Row 132
$ws.Range("H132").Value = 6179.3486
$ws.Range("I132").Value = 5895
$ws.Range("J132").Value = 6506.35
$ws.Range("K132").Value = 17685
$ws.Range("L132").Value = 19519.05
$ws.Range("M132").Value = -15155
$ws.Range("N132").Value = -24579.05
# Row 135
$ws.Range("H135").Value = 70087.22
$ws.Range("J135").Value = 70087.22
$ws.Range("L135").Value = 70087.22
$ws.Range("N135").Value = -80227.22
# Row 136
$ws.Range("H136").Value = 5547.2905
$ws.Range("I136").Value = 2732.5833
$ws.Range("J136").Value = 7325
$ws.Range("K136").Value = 8197.749899999999
$ws.Range("L136").Value = 21975
$ws.Range("M136").Value = -5647.749899999999
$ws.Range("N136").Value = -27075

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 28599808
# Row 67
$ws.Range("H67").Value = 28599808
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 82
$ws.Range("H82").Value = 39995
$ws.Range("I82").Value = 39995
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 39995
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -39612
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 39995
$ws.Range("I85").Value = 39995
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 39995
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -38669
$ws.Range("N85").ClearContents()
# Row 103
$ws.Range("H103").Value = 28928.143
$ws.Range("J103").Value = 28928.143
$ws.Range("L103").Value = 28928.143
$ws.Range("N103").Value = -31272.143
# Row 113
$ws.Range("H113").Value = 511.34482
$ws.Range("J113").Value = 534.8889
$ws.Range("L113").Value = 1604.6667
$ws.Range("N113").Value = -5944.6667
# Row 132
$ws.Range("H132").Value = 2238.8547
$ws.Range("I132").Value = 1547.3405
$ws.Range("J132").Value = 4405.6
$ws.Range("K132").Value = 4642.0215
$ws.Range("L132").Value = 13216.8
$ws.Range("M132").Value = -2112.0215
$ws.Range("N132").Value = -18276.8
# Row 136
$ws.Range("H136").Value = 2622
$ws.Range("I136").Value = 2231.4187
$ws.Range("K136").Value = 6694.256100000001
$ws.Range("M136").Value = -4144.256100000001
